# Update attendance summary values from 0 to 1 for the specific cells
# that represent the recorded attendance data (tut06 final code with comment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("H5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

$ws.Range("H12").Value = 1

$ws.Range("H13").Value = 1

$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1

$ws.Range("H15").Value = 1

$ws.Range("H16").Value = 1

$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 1
